# Actualización automática 2025-09-22 08:22:24
# Registra una venta de SAL SOLUBLE por 831.18 para el cliente
# "GONZALEZ CARDENAS ERNESTO PAOLO" (asesor LOZANO MOLINA TITO) en
# septiembre, y recalcula los totales dependientes en las tres hojas.

$wb = $excel.ActiveWorkbook

$wsVentasGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO: columna O = SAL SOLUBLE, fila 15 = cliente ---
$wsVentasGrupo.Range("O15").Value = 831.1799999999999
$wsVentasGrupo.Range("O32").Value = "1 de 30"

# --- VENTA MENSUAL: columna F = septiembre, fila 15 = cliente ---
$wsVentaMensual.Range("F15").Value = 831.1799999999999
$wsVentaMensual.Range("F32").Value = 9589.41

# --- CUMPLIMIENTO MENSUAL: fila 14 = SAL SOLUBLE, fila 15 = TOTAL ---
$wsCumplimiento.Range("D14").Value = 831.1799999999999
$wsCumplimiento.Range("E14").Value = -163.496851612446
$wsCumplimiento.Range("F14").Value = 1.244871915679299

$wsCumplimiento.Range("D15").Value = 9663.700000000001
$wsCumplimiento.Range("E15").Value = 22044.05990313501
$wsCumplimiento.Range("F15").Value = 0.3047739742423284

# El ancho de la columna E (5) en CUMPLIMIENTO MENSUAL creció por el
# nuevo valor negativo más largo (objetivo: ancho almacenado = 23).
$wsCumplimiento.Columns.Item(5).ColumnWidth = 22.166666666666668
